# Students can now select their degree type. Will ask them when they log in
# if it isn't set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the four students (shared-string text swaps).
$ws.Range("A2").Value2 = "Concepcion Hammes"
$ws.Range("A3").Value2 = "Pierce Rempel"
$ws.Range("A4").Value2 = "Estella Rogahn"
$ws.Range("A5").Value2 = "Jamil Thompson"

# Row 3 (Pierce Rempel) no longer applied/accepted; row 4 (Estella Rogahn) now is.
$ws.Range("C3:E3").Value2 = 0
$ws.Range("C4:E4").Value2 = 1

# Column A grew a little wider to fit the new (longer) names.
$ws.Columns.Item(1).ColumnWidth = 20.33
